$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task name updated: GPS Road Deviation -> GPS Road / Intersection Deviation
$ws.Range("A4").Value = "GPS Road / Intersection Deviation"

# Task marked complete: Progress Completion 0% -> 100%
$ws.Range("E4").Value = 1

# Reflect the user's active cell after editing the completion value
$ws.Range("E4").Select()
